$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above the current row 1094 (shifts 1094:1177 -> 1097:1180)
$ws.Rows("1094:1096").Insert()

# New row 1094 - Femacal de La Calera, Coquimbo, Platano, Sin especificar, Maduro
$ws.Range("A1094").Value = 3
$ws.Range("B1094").Value = "Femacal de La Calera"
$ws.Range("C1094").Value = "Coquimbo"
$ws.Range("D1094").Value = 44931
$ws.Range("E1094").Value = 5
$ws.Range("F1094").Value = "Fruta"
$ws.Range("G1094").Value = 100108
$ws.Range("H1094").Value = "Tropicales y subtropicales"
$ws.Range("I1094").Value = 100108006
$ws.Range("J1094").Value = "Plátano"
$ws.Range("K1094").Value = "Sin especificar"
$ws.Range("L1094").Value = "Maduro"
$ws.Range("M1094").Value = 120
$ws.Range("N1094").Value = 17000
$ws.Range("O1094").Value = 17000
$ws.Range("P1094").Value = 17000
$ws.Range("Q1094").Value = "`$/caja 20 kilos"
$ws.Range("R1094").Value = "Ecuador"
$ws.Range("S1094").Value = 850
$ws.Range("T1094").Value = 20

# New row 1095 - Pintón
$ws.Range("A1095").Value = 3
$ws.Range("B1095").Value = "Femacal de La Calera"
$ws.Range("C1095").Value = "Coquimbo"
$ws.Range("D1095").Value = 44931
$ws.Range("E1095").Value = 5
$ws.Range("F1095").Value = "Fruta"
$ws.Range("G1095").Value = 100108
$ws.Range("H1095").Value = "Tropicales y subtropicales"
$ws.Range("I1095").Value = 100108006
$ws.Range("J1095").Value = "Plátano"
$ws.Range("K1095").Value = "Sin especificar"
$ws.Range("L1095").Value = "Pintón"
$ws.Range("M1095").Value = 200
$ws.Range("N1095").Value = 18000
$ws.Range("O1095").Value = 18000
$ws.Range("P1095").Value = 18000
$ws.Range("Q1095").Value = "`$/caja 20 kilos"
$ws.Range("R1095").Value = "Ecuador"
$ws.Range("S1095").Value = 900
$ws.Range("T1095").Value = 20

# New row 1096 - Primera Pintón
$ws.Range("A1096").Value = 3
$ws.Range("B1096").Value = "Femacal de La Calera"
$ws.Range("C1096").Value = "Coquimbo"
$ws.Range("D1096").Value = 44931
$ws.Range("E1096").Value = 5
$ws.Range("F1096").Value = "Fruta"
$ws.Range("G1096").Value = 100108
$ws.Range("H1096").Value = "Tropicales y subtropicales"
$ws.Range("I1096").Value = 100108006
$ws.Range("J1096").Value = "Plátano"
$ws.Range("K1096").Value = "Sin especificar"
$ws.Range("L1096").Value = "Primera Pintón"
$ws.Range("M1096").Value = 360
$ws.Range("N1096").Value = 20000
$ws.Range("O1096").Value = 21000
$ws.Range("P1096").Value = 20556
$ws.Range("Q1096").Value = "`$/caja 20 kilos"
$ws.Range("R1096").Value = "Ecuador"
$ws.Range("S1096").Value = 1028
$ws.Range("T1096").Value = 20
